{"js": "const pairs = [[\"2025-07-26 Saturday\", \"2025-07-27 Sunday\"], [\"565\u00d72=1130\", \"391\u00d72=782\"], [\"702\u00d74=2808\", \"462\u00d72=924\"], [\"380\u00d73=1140\", \"259\u00d76=1554\"], [\"256\u00d72=512\", \"674\u00d74=2696\"], [\"528\u00d75=2640\", \"631\u00d73=1893\"], [\"730\u00d76=4380\", \"735\u00d72=1470\"], [\"407\u00d78=3256\", \"904\u00d75=4520\"], [\"116\u00d77=812\", \"295\u00d73=885\"], [\"428\u00d79=3852\", \"762\u00d75=3810\"], [\"354\u00d72=708\", \"122\u00d73=366\"], [\"204\u00d76=1224\", \"613\u00d73=1839\"], [\"693\u00d78=5544\", \"733\u00d75=3665\"], [\"140\u00d72=280\", \"113\u00d73=339\"], [\"330\u00d72=660\", \"710\u00d75=3550\"], [\"601\u00d73=1803\", \"188\u00d74=752\"], [\"385\u00d79=3465\", \"586\u00d74=2344\"], [\"665\u00d79=5985\", \"138\u00d72=276\"], [\"599\u00d76=3594\", \"780\u00d76=4680\"], [\"246\u00d73=738\", \"945\u00d78=7560\"], [\"325\u00d74=1300\", \"204\u00d79=1836\"], [\"686\u00d75=3430\", \"738\u00d74=2952\"], [\"333\u00d72=666\", \"443\u00d74=1772\"], [\"655\u00d75=3275\", \"565\u00d73=1695\"], [\"663\u00d77=4641\", \"806\u00d72=1612\"], [\"498\u00d76=2988\", \"209\u00d73=627\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-26 Saturday\", \"2025-07-27 Sunday\"),\n    @(\"565\u00d72=1130\", \"391\u00d72=782\"),\n    @(\"702\u00d74=2808\", \"462\u00d72=924\"),\n    @(\"380\u00d73=1140\", \"259\u00d76=1554\"),\n    @(\"256\u00d72=512\", \"674\u00d74=2696\"),\n    @(\"528\u00d75=2640\", \"631\u00d73=1893\"),\n    @(\"730\u00d76=4380\", \"735\u00d72=1470\"),\n    @(\"407\u00d78=3256\", \"904\u00d75=4520\"),\n    @(\"116\u00d77=812\", \"295\u00d73=885\"),\n    @(\"428\u00d79=3852\", \"762\u00d75=3810\"),\n    @(\"354\u00d72=708\", \"122\u00d73=366\"),\n    @(\"204\u00d76=1224\", \"613\u00d73=1839\"),\n    @(\"693\u00d78=5544\", \"733\u00d75=3665\"),\n    @(\"140\u00d72=280\", \"113\u00d73=339\"),\n    @(\"330\u00d72=660\", \"710\u00d75=3550\"),\n    @(\"601\u00d73=1803\", \"188\u00d74=752\"),\n    @(\"385\u00d79=3465\", \"586\u00d74=2344\"),\n    @(\"665\u00d79=5985\", \"138\u00d72=276\"),\n    @(\"599\u00d76=3594\", \"780\u00d76=4680\"),\n    @(\"246\u00d73=738\", \"945\u00d78=7560\"),\n    @(\"325\u00d74=1300\", \"204\u00d79=1836\"),\n    @(\"686\u00d75=3430\", \"738\u00d74=2952\"),\n    @(\"333\u00d72=666\", \"443\u00d74=1772\"),\n    @(\"655\u00d75=3275\", \"565\u00d73=1695\"),\n    @(\"663\u00d77=4641\", \"806\u00d72=1612\"),\n    @(\"498\u00d76=2988\", \"209\u00d73=627\")\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
